$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INFORME DICIEMBRE")

# H13 and H15: mark as reviewed by CARLOS MEZA (same shared string already used in H11, H29, H31, H33, H35)
$ws.Range("H13").Value = "CARLOS MEZA"
$ws.Range("H15").Value = "CARLOS MEZA"

# Row 16: fill in the week's dates (continuing the sequence from row 14)
$ws.Range("A16").Value = 44186
$ws.Range("B16").Value = 44187
$ws.Range("C16").Value = 44188
$ws.Range("D16").Value = 44189
$ws.Range("E16").Value = 44190
$ws.Range("F16").Value = 44191
$ws.Range("G16").Value = 44192

# Row 17: add the week's notes
$ws.Range("B17").Value = "se organizo el menu principal, y se configuro el acceso por Roles de Usuario"
$ws.Range("C17").Value = "Modificaciones en Formularios pendientes"
$ws.Range("D17").Value = "-"
$ws.Range("E17").Value = "Modificaciones en formulario(Banco)"

$ws.Range("H17").Value = "CARLOS MEZA"
$ws.Range("H17").HorizontalAlignment = -4108
$ws.Rows.Item(17).RowHeight = 71.25

$ws.Range("E17").Select()
